# The workbook tracks weekly fruit/vegetable prices. A new weekly record was
# added as row 129 (same Mercado / Producto / Variedad / Calidad as the
# previous row 129, but a newer date and a different Volumen), which pushed
# every subsequent data row down by one (old row 129 -> new row 130, ...,
# old row 195 -> new row 196).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129; Excel shifts rows 129-195 down to 130-196 and
# copies the formatting (including the date style on column D) from the
# row above.
$ws.Rows("129:129").Insert()

# Populate the newly inserted row 129 with the new weekly record.
$ws.Range("A129").Value = 1
$ws.Range("B129").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C129").Value = "Arica y Parinacota"
$ws.Range("D129").Value = 44603
$ws.Range("E129").Value = 15
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100102
$ws.Range("H129").Value = "Cítricos"
$ws.Range("I129").Value = 100102003
$ws.Range("J129").Value = "Limón"
$ws.Range("K129").Value = "Tahití"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 350
$ws.Range("N129").Value = 35000
$ws.Range("O129").Value = 36000
$ws.Range("P129").Value = 35500
$ws.Range("Q129").Value = "$/caja 24 kilos"
$ws.Range("R129").Value = "Perú"
$ws.Range("S129").Value = 1479
$ws.Range("T129").Value = 24

$wb.Save()
